$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update header row (row 1) ---
$ws.Range("A1").Value = "species"
$ws.Range("B1").Value = "upper_CI"
$ws.Range("C1").Value = "lower_CI"
$ws.Range("D1").Value = "mean_AIC"
$ws.Range("E1").Value = "iterations"
$ws.Range("F1").Value = "tgam_AIC"

# --- Update species labels in column A (rows 2-9) ---
$ws.Range("A2").Value = "arrowtooth"
$ws.Range("A3").Value = "english"
$ws.Range("A4").Value = "sanddab"
$ws.Range("A5").Value = "dover"
$ws.Range("A6").Value = "rex"
$ws.Range("A7").Value = "lingcod"
$ws.Range("A8").Value = "petrale"
$ws.Range("A9").Value = "sablefish"

# --- Update the active selection ---
$ws.Range("A11").Select() | Out-Null
